# Atualização de bases das ligas, do dia: 09-03-2024 às 13:07
# Croatia HNL sheet update:
#  - The match previously at row 127 (id=125, NK Varazdin vs NK Rudes, B=6769304)
#    is removed from the feed; subsequent rows shift up by one.
#  - Row 126 (id=124, Dinamo Zagreb vs HNK Gorica) has now been played: final
#    score / result columns (H, I, J) are filled in and its closing odds /
#    Asian-handicap P&L columns are refreshed.
#  - The rows that shifted up (old 128/129/130 -> new 127/128/129) get
#    refreshed closing-odds values (R, S, U, V and, for the first of them,
#    O, P too) reflecting the newer scrape.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the old row 127 (id=125) entirely -- this naturally shifts
#    rows 128..130 up to 127..129, matching the diff's row-129/130 removal.
$ws.Rows(127).Delete()

# The engine's row-shift implementation auto-extends column A's numeric
# series (an AutoFill-like heuristic) when rows move up; restore the
# plain "id" values the shifted rows actually carry (unchanged by the
# deletion -- id is just data, not a formula/series).
$ws.Cells.Item(126, 1).Value = 124
$ws.Cells.Item(127, 1).Value = 125
$ws.Cells.Item(128, 1).Value = 126
$ws.Cells.Item(129, 1).Value = 127

# 2) Row 126 (id=124): match has been played now -- add result + refresh
#    closing odds / PL columns.
$ws.Cells.Item(126, 8).Value  = 2          # H126 FTHG
$ws.Cells.Item(126, 9).Value  = 2          # I126 FTAG
$ws.Cells.Item(126, 10).Value = "D"        # J126 FTR
$ws.Cells.Item(126, 15).Value = 2.9        # O126 oddD
$ws.Cells.Item(126, 16).Value = 3.3        # P126 oddA
$ws.Cells.Item(126, 18).Value = 1.925      # R126 oddAHH
$ws.Cells.Item(126, 19).Value = 1.925      # S126 oddAHA
$ws.Cells.Item(126, 20).Value = 1.75       # T126 AhOU
$ws.Cells.Item(126, 21).Value = 1.775      # U126 oddAHOver
$ws.Cells.Item(126, 22).Value = 2.1        # V126 oddAHUnder
$ws.Cells.Item(126, 23).Value = -1         # W126 PLH
$ws.Cells.Item(126, 24).Value = 1.9        # X126 PLD
$ws.Cells.Item(126, 25).Value = -1         # Y126 PLA
$ws.Cells.Item(126, 26).Value = -0.5       # Z126 PL_Ahh
$ws.Cells.Item(126, 27).Value = 0.4625     # AA126 PL_Aha
$ws.Cells.Item(126, 28).Value = 0.7749999999999999  # AB126 PL_AhOver
$ws.Cells.Item(126, 29).Value = -1         # AC126 PL_AhUnder

# 3) Row 127 (after the shift; id=125, HNK Gorica vs Hajduk Split): refresh
#    closing-odds columns with the newer scrape values.
$ws.Cells.Item(127, 15).Value = 3.6        # O127 oddD
$ws.Cells.Item(127, 16).Value = 1.5        # P127 oddA
$ws.Cells.Item(127, 18).Value = 1.975      # R127 oddAHH
$ws.Cells.Item(127, 19).Value = 1.875      # S127 oddAHA
$ws.Cells.Item(127, 21).Value = 1.95       # U127 oddAHOver
$ws.Cells.Item(127, 22).Value = 1.9        # V127 oddAHUnder

# 4) Row 128 (after the shift; id=126, HNK Rijeka vs NK Osijek): refresh
#    closing-odds columns.
$ws.Cells.Item(128, 18).Value = 1.95       # R128 oddAHH
$ws.Cells.Item(128, 19).Value = 1.9        # S128 oddAHA
$ws.Cells.Item(128, 21).Value = 1.9        # U128 oddAHOver
$ws.Cells.Item(128, 22).Value = 1.95       # V128 oddAHUnder

# 5) Row 129 (after the shift; id=127, Dinamo Zagreb vs Slaven Belupo):
#    refresh closing-odds columns.
$ws.Cells.Item(129, 18).Value = 1.825      # R129 oddAHH
$ws.Cells.Item(129, 19).Value = 2.025      # S129 oddAHA
$ws.Cells.Item(129, 21).Value = 1.95       # U129 oddAHOver
$ws.Cells.Item(129, 22).Value = 1.9        # V129 oddAHUnder
